$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Copy the O column's current number format (hh:mm, style index 10) into the
#        new P cells for rows 62-69 before we touch any values, so the new P cells
#        inherit the same style as O (bold hh:mm) rather than the "Frecuencia_Min"
#        style used elsewhere (s=12).
$ws.Range("O62:O69").Copy()
$ws.Range("P62:P69").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 2) Re-point rows 62-69 from the "Flexible" schedule type to "Frecuencia"
#        (frequency-based). Typing the literal text lets Excel reuse/garbage
#        collect the shared-string table once "Flexible" has no more referrers.
$ws.Cells.Item(62, 14).Value = "Frecuencia"
$ws.Cells.Item(63, 14).Value = "Frecuencia"
$ws.Cells.Item(64, 14).Value = "Frecuencia"
$ws.Cells.Item(65, 14).Value = "Frecuencia"
$ws.Cells.Item(66, 14).Value = "Frecuencia"
$ws.Cells.Item(67, 14).Value = "Frecuencia"
$ws.Cells.Item(68, 14).Value = "Frecuencia"
$ws.Cells.Item(69, 14).Value = "Frecuencia"

# --- 3) The old single "duration" value in column O becomes the new
#        "Frecuencia_Min" value in column P; O is reset to 0 (no longer used
#        once the row is frequency-based).
$ws.Range("P62").Value = 0.010416666666666666
$ws.Range("P63").Value = 0.010416666666666666
$ws.Range("P64").Value = 0.0083333333333333332
$ws.Range("P65").Value = 0.0083333333333333332
$ws.Range("P66").Value = 0.017361111111111112
$ws.Range("P67").Value = 0.017361111111111112
$ws.Range("P68").Value = 0.024305555555555556
$ws.Range("P69").Value = 0.024305555555555556

$ws.Range("O62:O69").Value = 0

# --- 4) Widen column P (Frecuencia_Min) now that it carries real data for more
#        rows, and drop the old "best fit" auto-sizing in favour of a fixed width.
$ws.Columns.Item(16).ColumnWidth = 23.74

# --- 5) Leave the selection where the user finished editing.
$ws.Range("O63").Select() | Out-Null
